# Turkish "Neden Epic - BTC karsilastirmasi" deck:
# Fix the typo/wording in the "DEGISTIRILEBILIRLIK" label textbox so it
# reads "DEGISTIRILEBILIR" (drop the trailing "LIK"), matching the other
# single-word feature labels on the slide (e.g. "OLCEKLENEBILIR").

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$found = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    if ($shp.Name -eq "TextBox 87") {
        $found = $shp
        break
    }
}

if ($found -ne $null) {
    $found.TextFrame.TextRange.Text = "DEĞİŞTİRİLEBİLİR"
}
